# Applies the cryptos.xlsx price/volume refresh described in the commit
# "Updated cryptos list on Wed Oct  4 17:50:49 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.560.27"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "1.639.60"
$ws.Range("E3").Value = "  -0.96%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.537"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.81%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.03"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.74%  "
$ws.Range("E9").Value = "  -1.82%  "
$ws.Range("E11").Value = "  +1.41%  "
$ws.Range("D12").Value = "1.872.07"
$ws.Range("D13").Value = "1.639.03"
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.562"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.09"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.72%  "
$ws.Range("D17").Value = "27.541.08"
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.92%  "
$ws.Range("D20").Value = "0.0₃0723"
$ws.Range("E20").Value = "  -0.43%  "
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.22%  "
$ws.Range("E24").Value = "  -3.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.64%  "
$ws.Range("E27").Value = "  +1.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.59"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.23%  "
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("E31").Value = "  -2.12%  "
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.15"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.35%  "
$ws.Range("D34").Value = "1.426.25"
$ws.Range("E34").Value = "  -2.53%  "
$ws.Range("E35").Value = "  +1.77%  "
$ws.Range("E36").Value = "  -1.96%  "
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.877"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.77%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0166"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.01%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.907"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +16.10%  "
$ws.Range("E41").Value = "  -1.61%  "
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("E44").Value = "  +0.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.83"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.07%  "
$ws.Range("D47").Value = "1.780.70"
$ws.Range("E47").Value = "  -0.93%  "
$ws.Range("E48").Value = "  -2.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.26"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.45%  "
$ws.Range("D50").Value = "0.0₆0106"
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("E51").Value = "  -2.70%  "
